# Apply sprint-task sheet updates:
#  - Mark all S10 tasks (rows 71-82) as "implemented"
#  - Append new S11 tasks (rows 83-88) describing the Data Manager / OHLCV cache work

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Flip status column (G) for the existing S10 rows (71-82) to "implemented"
# ---------------------------------------------------------------------------
for ($r = 71; $r -le 82; $r++) {
    $ws.Cells.Item($r, 7).Value = "implemented"
}

# ---------------------------------------------------------------------------
# 2) Append the new S11 rows (83-88)
# ---------------------------------------------------------------------------
$newRows = @(
    ,@(
        "S11",
        "G01",
        "Data Manager & OHLCV cache: PRD and design",
        "S11_G01_TB001",
        "Draft PRD for persistent OHLCV cache and Data Manager, defining base timeframe, horizon, and ensure_coverage semantics.",
        "See docs/qlab_data_cache_prd.md for detailed design.",
        "implemented"
    )
    ,@(
        "S11",
        "G02",
        "Data Manager & OHLCV cache: backend implementation",
        "S11_G02_TB001",
        "Implement DataManager helpers to compute coverage gaps per symbol/timeframe using price_bars and price_fetches.",
        "Reuses existing prices DB schema; no breaking changes.",
        "pending"
    )
    ,@(
        "S11",
        "G02",
        "Data Manager & OHLCV cache: backend implementation",
        "S11_G02_TB002",
        "Wire DataManager.ensure_symbol_coverage into run_single_backtest and ensure_group_coverage into run_group_backtest.",
        "Backtests should no longer call Kite/yfinance directly; they rely on the local cache.",
        "pending"
    )
    ,@(
        "S11",
        "G02",
        "Data Manager & OHLCV cache: backend implementation",
        "S11_G02_TB003",
        "Add regression tests that run backtests without prior Fetch Data calls and assert that coverage is built and reused.",
        "Tests may use a synthetic provider or stub DataService to avoid real network calls.",
        "pending"
    )
    ,@(
        "S11",
        "G03",
        "Data Manager & OHLCV cache: Data page integration",
        "S11_G03_TF001",
        "Add a switch on the Data page to choose between casual preview and saving fetched data to the persistent cache.",
        "When saving, default timeframe to base_timeframe and extend duration to the configured BT horizon.",
        "pending"
    )
    ,@(
        "S11",
        "G03",
        "Data Manager & OHLCV cache: Data page integration",
        "S11_G03_TF002",
        "Update Coverage Summary UI to indicate which rows are BT-ready cache entries vs preview-only data.",
        "Leverages existing created_at and coverage_id fields; may add a simple badge/flag.",
        "pending"
    )
)

$r = 83
foreach ($rowData in $newRows) {
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowData[$c]
    }
    $r = $r + 1
}
